# Update Excel files after daily scrape - 2025-12-04 03:29:16 UTC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-16), columns A-H
$data = @(
    @("1330354", "https://aiesec.org/opportunity/global-talent/1330354", "SCM & BI  Senior Intern", "Panamá, Provincia de Panamá, Panamá", "No", "8 applicants", "6 - 18 Months", "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"),
    @("1330335", "https://aiesec.org/opportunity/global-talent/1330335", "Guest Relations Executive", "Weligama, Sri Lanka", "No", "0 applicants", "3 - 6 Months", "Weligama Cliff"),
    @("1330306", "https://aiesec.org/opportunity/global-talent/1330306", "Behavioral Marketing & Communications Intern", "Amman, Jordan", "No", "5 applicants", "9 - 12 Weeks", "Amoux Group"),
    @("1330301", "https://aiesec.org/opportunity/global-talent/1330301", "Business Development Officer", "Amman, Jordan", "No", "8 applicants", "9 - 12 Weeks", "International TEFL Training Institute"),
    @("1330293", "https://aiesec.org/opportunity/global-talent/1330293", "Guest Relations Manager", "Mumbai, Maharashtra, India", "No", "2 applicants", "6 - 18 Months", "Taj Mahal Hotels & Tower"),
    @("1330246", "https://aiesec.org/opportunity/global-talent/1330246", "Mobile Developer", "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Safarni"),
    @("1330244", "https://aiesec.org/opportunity/global-talent/1330244", "Web developer", "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Safarni"),
    @("1330236", "https://aiesec.org/opportunity/global-talent/1330236", "IT Sales", "Hanoi, Vietnam", "No", "5 applicants", "9 - 12 Weeks", "Vitex Vietnam Software Joint Stock Company"),
    @("1329876", "https://aiesec.org/opportunity/global-talent/1329876", "[Impact Porto Alegre] - Import and Product Development Intern (Mandarin Chinese Speaker)", "Porto Alegre, RS, Brazil", "No", "13 applicants", "6 - 18 Months", "Asialink - International Commerce"),
    @("1329104", "https://aiesec.org/opportunity/global-talent/1329104", "Sales Intern", "Pannipitiya, Sri Lanka", "No", "9 applicants", "6 - 18 Months", "Frella International"),
    @("1328532", "https://aiesec.org/opportunity/global-talent/1328532", "Marketing Strategist & Media Buyer", "Cairo, Cairo Governorate, Egypt", "No", "12 applicants", "3 - 6 Months", "Al Mowared for Engineering works"),
    @("1327503", "https://aiesec.org/opportunity/global-talent/1327503", "Taste Hungary | Customer Service trainee", "Budapest, Magyarország", "No", "44 applicants", "3 - 6 Months", "JO&JOE Budapest"),
    @("1327149", "https://aiesec.org/opportunity/global-talent/1327149", "Data Engineering Intern", "Nugegoda, Sri Lanka", "No", "73 applicants", "3 - 6 Months", "Altria Consulting (PVT) LTD"),
    @("1323077", "https://aiesec.org/opportunity/global-talent/1323077", "Pharmacy Intern", "Manipal, Karnataka, India", "No", "18 applicants", "9 - 12 Weeks", "M.A.H.E."),
    @("1307153", "https://aiesec.org/opportunity/global-talent/1307153", "Accelerate Romania | Technology Business Development Consultant", "Bucharest, Romania", "No", "158 applicants", "6 - 18 Months", "Rinf Tech")
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        if ($col -eq 1) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$col - 1]
        if ($col -eq 1 -or $col -eq 5) {
            $cell.ClearFormats()
        }
    }
    $rowIndex++
}

# Update column widths (COM ColumnWidth adds ~0.8333 vs the raw OOXML
# "width" attribute when saved; subtract that offset so the serialized
# width matches the target exactly)
$ws.Columns.Item(3).ColumnWidth = 91 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 57 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 17 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 60 - (5/6)
